$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-10 (columns A-T), reflecting the updated TPM-based NATMI
# ligand-receptor (Efnb3-Rhbdl2) results. Row 2 corresponds to sheet row 2, etc.
$data = @(
    @("ECs","Efnb3","Rhbdl2","ECs",2,0.6666666666666666,0.2151756666666667,0.645527,0.1791915537270423,0.1791915537270423,3,1,2.214002666666667,6.642008000000001,0.9692171064132561,0.9692171064132561,0.4763994998017778,4.287595498216001,0.1736755191970195,0.1736755191970195),
    @("ECs","Efnb3","Rhbdl2","MuSCs",2,0.6666666666666666,0.2151756666666667,0.645527,0.1791915537270423,0.1791915537270423,1,0.3333333333333333,0.06541966666666667,0.196259,0.02863856533860833,0.02863856533860832,0.01407672038811111,0.126690483493,0.005131789019538645,0.005131789019538644),
    @("ECs","Efnb3","Rhbdl2","Resolving-Mac",2,0.6666666666666666,0.2151756666666667,0.645527,0.1791915537270423,0.1791915537270423,1,0.3333333333333333,0.004898333333333334,0.014695,0.002144328248135624,0.002144328248135624,0.001054002140555556,0.009486019264999999,0.0003842455104842091,0.0003842455104842091),
    @("FAPs","Efnb3","Rhbdl2","ECs",2,0.6666666666666666,0.5010686666666667,1.503206,0.4172742870736815,0.4172742870736815,3,1,2.214002666666667,6.642008000000001,0.9692171064132561,0.9692171064132561,1.109367364183111,9.984306277648001,0.4044293770982079,0.4044293770982079),
    @("FAPs","Efnb3","Rhbdl2","MuSCs",2,0.6666666666666666,0.5010686666666667,1.503206,0.4172742870736815,0.4172742870736815,1,0.3333333333333333,0.06541966666666667,0.196259,0.02863856533860833,0.02863856533860832,0.03277974515044445,0.295017706354,0.01195013693448083,0.01195013693448083),
    @("FAPs","Efnb3","Rhbdl2","Resolving-Mac",2,0.6666666666666666,0.5010686666666667,1.503206,0.4172742870736815,0.4172742870736815,1,0.3333333333333333,0.004898333333333334,0.014695,0.002144328248135624,0.002144328248135624,0.002454401352222222,0.02208961217,0.0008947730409927487,0.0008947730409927487),
    @("MuSCs","Efnb3","Rhbdl2","ECs",2,0.6666666666666666,0.4845693333333334,1.453708,0.4035341591992763,0.4035341591992763,3,1,2.214002666666667,6.642008000000001,0.9692171064132561,0.9692171064132561,1.072837796184889,9.655540165664,0.3911122101180288,0.3911122101180288),
    @("MuSCs","Efnb3","Rhbdl2","MuSCs",2,0.6666666666666666,0.4845693333333334,1.453708,0.4035341591992763,0.4035341591992763,1,0.3333333333333333,0.06541966666666667,0.196259,0.02863856533860833,0.02863856533860832,0.03170036426355555,0.285303278372,0.01155663938458885,0.01155663938458885),
    @("MuSCs","Efnb3","Rhbdl2","Resolving-Mac",2,0.6666666666666666,0.4845693333333334,1.453708,0.4035341591992763,0.4035341591992763,1,0.3333333333333333,0.004898333333333334,0.014695,0.002144328248135624,0.002144328248135624,0.002373582117777778,0.02136223906,0.000865309696658666,0.000865309696658666)
)

$numRows = $data.Count
$numCols = 20

for ($r = 0; $r -lt $numRows; $r++) {
    for ($c = 0; $c -lt $numCols; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $data[$r][$c]
    }
}
